$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B..I are always the constant "N/A" placeholder (shared string index 11 in the original file)
$na = "N/A"

# Row 2: Sequence 1
$ws.Range("A2").Value = '2025-11-24T09:31:23.964Z'
$ws.Range("B2:I2").Value = $na
$ws.Range("J2").Value = '{"Sequence":1,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"d675a633-726c-4afd-b83e-b772176abf33","EventDtm":"2025-11-24T09:31:03Z","AppDtm":"2025-11-24T09:31:23Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-24T09:31:03Z","DeviceID":"JSGA623040320","GPSLockState":"LOCKED","GPSSatelliteCount":15,"GPSLastLock":0,"GPSLatitude":17.657838,"GPSLongitude":83.101784,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"49","LAC":"31121","CID":"248666913","RSSI":"-73","ExtPower":true,"ExtPowerVoltage":29.7,"BatteryVoltage":7.9,"DeviceTemp":36,"RTDLOn":false,"VerFW":"W0206 1.91","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"KKFU6712079","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-24T09:31:03Z","AssetType":"Reefer","AssetID":"KKFU6712079","OEM":"CARRIER","TAmb":29.85,"TAmbQ":null,"TUSDA4":-50,"TUSDA4Q":"OOR","PctCO2":0,"PctCO2Q":null,"PctCO2Set":0,"PctCO2SetQ":null,"PSuc":65.86,"PSucQ":null,"TDis":-191.8,"TDisQ":"OOR","FreqComp":null,"TSuc":0,"TSucQ":null,"MCond":0,"PCond":137.31,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":"Idle","SnCtrl":"04742313","AmpPhA":0.92,"AmpPhB":1.54,"AmpPhC":0.92,"PDis":150.14,"PDisQ":null,"PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":0,"PctExpVlvQ":null,"TEvap":20.4,"TEvapQ":null,"MCtrl3":null,"PctHtr":null,"PctHtrQ":null,"MEvapFanHS":0,"PctGasVlv":null,"PctGasVlvQ":null,"PctHum":0,"PctHumQ":null,"PctHumSet":0,"PctHumSetQ":"OOR","FreqLine":50,"FreqLineQ":null,"VLine1":412.26,"VLine2":null,"VLine3":null,"MEvapFanLS":1,"PctO2":0,"PctO2Q":"asProvided","PctO2Set":0,"PctO2SetQ":"OOR","MCtrl2":null,"TRtn1":20.09,"TRtn1Q":null,"TRtn2":20.1,"TRtn2Q":null,"TSet":20,"TSetQ":null,"VerSWMajor":"5154","VerSWMinor":null,"PctSucVlv":18.49,"PctSucVlvQ":null,"TSup1":18.99,"TSup1Q":null,"TSup2":19.07,"TSup2Q":null,"AmpTotalDraw":null,"AmpTotalDrawQ":null,"TUSDA1":-50,"TUSDA1Q":"OOR","TUSDA2":-50,"TUSDA2Q":"OOR","TUSDA3":-50,"TUSDA3Q":"OOR","CmhVent":0,"CmhVentQ":null,"BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2024-03-28T10:32:52Z","PTIResult":"Passed","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":null,"ReeferAlarms":null}}}'

# Row 3: Sequence 2
$ws.Range("A3").Value = '2025-11-24T09:31:36.082Z'
$ws.Range("B3:I3").Value = $na
$ws.Range("J3").Value = '{"Sequence":2,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"92c925a0-c632-4961-ba83-f560fef5c6e6","EventDtm":"2025-11-24T09:31:20Z","AppDtm":"2025-11-24T09:31:34Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-24T09:31:20Z","DeviceID":"JSGA623040277","GPSLockState":"LOCKED","GPSSatelliteCount":13,"GPSLastLock":0,"GPSLatitude":17.572157,"GPSLongitude":78.514977,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"49","LAC":"19327","CID":"235323414","RSSI":"-67","ExtPower":true,"ExtPowerVoltage":33.1,"BatteryVoltage":8,"DeviceTemp":33,"RTDLOn":false,"VerFW":"W0206 1.91","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"DFOU6120281","LastAssetRunState":"Running"},"ReeferData":null}}'

# Row 4: Sequence 3
$ws.Range("A4").Value = '2025-11-24T09:31:41.777Z'
$ws.Range("B4:I4").Value = $na
$ws.Range("J4").Value = '{"Sequence":3,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"96ac87a5-715d-4072-91ee-515eb562ea66","EventDtm":"2025-11-24T09:31:24Z","AppDtm":"2025-11-24T09:31:41Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-24T09:31:24Z","DeviceID":"JSGA623040262","GPSLockState":"LOCKED","GPSSatelliteCount":18,"GPSLastLock":0,"GPSLatitude":12.805008,"GPSLongitude":77.662269,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"45","LAC":"9003","CID":"46904076","RSSI":null,"ExtPower":true,"ExtPowerVoltage":30,"BatteryVoltage":8,"DeviceTemp":34,"RTDLOn":false,"VerFW":"W0206 1.75","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"TRIU6681542","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-24T09:31:24Z","AssetType":"Reefer","AssetID":"TRIU6681542","OEM":"CARRIER","TAmb":28.98,"TAmbQ":null,"TUSDA4":-50,"TUSDA4Q":"OOR","PctCO2":5,"PctCO2Q":null,"PctCO2Set":1,"PctCO2SetQ":null,"PSuc":-14.73,"PSucQ":"OOR","TDis":-196.88,"TDisQ":"OOR","FreqComp":null,"TSuc":0,"TSucQ":null,"MCond":0,"PCond":168.39,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":"Cool","SnCtrl":"04958167","AmpPhA":10.08,"AmpPhB":10.52,"AmpPhC":10.08,"PDis":-14.6,"PDisQ":"OOR","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":0,"PctExpVlvQ":null,"TEvap":6.02,"TEvapQ":null,"MCtrl3":null,"PctHtr":null,"PctHtrQ":null,"MEvapFanHS":1,"PctGasVlv":null,"PctGasVlvQ":null,"PctHum":96.42,"PctHumQ":null,"PctHumSet":0,"PctHumSetQ":"OOR","FreqLine":50,"FreqLineQ":null,"VLine1":413.49,"VLine2":null,"VLine3":null,"MEvapFanLS":0,"PctO2":-0.05,"PctO2Q":"OOR","PctO2Set":3,"PctO2SetQ":null,"MCtrl2":null,"TRtn1":5.48,"TRtn1Q":null,"TRtn2":5.48,"TRtn2Q":null,"TSet":5,"TSetQ":null,"VerSWMajor":"5178","VerSWMinor":null,"PctSucVlv":3.17,"PctSucVlvQ":null,"TSup1":4.96,"TSup1Q":null,"TSup2":4.98,"TSup2Q":null,"AmpTotalDraw":null,"AmpTotalDrawQ":null,"TUSDA1":-50,"TUSDA1Q":"OOR","TUSDA2":-50,"TUSDA2Q":"OOR","TUSDA3":-50,"TUSDA3Q":"OOR","CmhVent":0,"CmhVentQ":null,"BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2025-05-12T09:58:17Z","PTIResult":"Skipped","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":null,"ReeferAlarms":null}}}'

# Row 5: Sequence 4
$ws.Range("A5").Value = '2025-11-24T09:32:42.721Z'
$ws.Range("B5:I5").Value = $na
$ws.Range("J5").Value = '{"Sequence":4,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"991ea4d3-36a7-491a-996d-3ecdf304092d","EventDtm":"2025-11-24T09:32:26Z","AppDtm":"2025-11-24T09:32:42Z","Events":["BatteryPowerOn"]},"DeviceData":{"DeviceDataDtm":"2025-11-24T09:32:26Z","DeviceID":"JSGA623040311","GPSLockState":"LOCKED","GPSSatelliteCount":13,"GPSLastLock":0,"GPSLatitude":26.606171,"GPSLongitude":80.723755,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"15","LAC":"55363","CID":"210193709","RSSI":null,"ExtPower":false,"ExtPowerVoltage":6,"BatteryVoltage":8,"DeviceTemp":36,"RTDLOn":false,"VerFW":"W0206 1.75","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":null,"LastAssetRunState":null},"ReeferData":null}}'

# Row 6: Sequence 5
$ws.Range("A6").Value = '2025-11-24T09:33:15.758Z'
$ws.Range("B6:I6").Value = $na
$ws.Range("J6").Value = '{"Sequence":5,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"6f9d257f-37a2-4637-bfef-217051800198","EventDtm":"2025-11-24T09:33:01Z","AppDtm":"2025-11-24T09:33:15Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-24T09:33:01Z","DeviceID":"JSGA622340343","GPSLockState":"LOCKED","GPSSatelliteCount":14,"GPSLastLock":0,"GPSLatitude":12.805058,"GPSLongitude":77.662234,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"45","LAC":"9003","CID":"46904076","RSSI":null,"ExtPower":true,"ExtPowerVoltage":29.9,"BatteryVoltage":8,"DeviceTemp":33,"RTDLOn":false,"VerFW":"W0206 1.75","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"CCLU1035976","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-24T09:33:01Z","AssetType":"Reefer","AssetID":"CCLU1035976","OEM":"CARRIER","TAmb":33.34,"TAmbQ":null,"TUSDA4":-50,"TUSDA4Q":"OOR","PctCO2":5,"PctCO2Q":null,"PctCO2Set":1,"PctCO2SetQ":null,"PSuc":-14.73,"PSucQ":"OOR","TDis":-196.88,"TDisQ":"OOR","FreqComp":null,"TSuc":0,"TSucQ":null,"MCond":1,"PCond":106.15,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":"Cool","SnCtrl":"04904656","AmpPhA":11.32,"AmpPhB":10.61,"AmpPhC":11.32,"PDis":-14.6,"PDisQ":"OOR","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":0,"PctExpVlvQ":null,"TEvap":6.26,"TEvapQ":null,"MCtrl3":null,"PctHtr":null,"PctHtrQ":null,"MEvapFanHS":1,"PctGasVlv":null,"PctGasVlvQ":null,"PctHum":0,"PctHumQ":null,"PctHumSet":0,"PctHumSetQ":"OOR","FreqLine":50,"FreqLineQ":null,"VLine1":414.82,"VLine2":null,"VLine3":null,"MEvapFanLS":0,"PctO2":0.1,"PctO2Q":"unknown","PctO2Set":3,"PctO2SetQ":null,"MCtrl2":null,"TRtn1":5.82,"TRtn1Q":null,"TRtn2":5.84,"TRtn2Q":null,"TSet":5,"TSetQ":null,"VerSWMajor":"5180","VerSWMinor":null,"PctSucVlv":5.37,"PctSucVlvQ":null,"TSup1":5.11,"TSup1Q":null,"TSup2":5.03,"TSup2Q":null,"AmpTotalDraw":null,"AmpTotalDrawQ":null,"TUSDA1":-50,"TUSDA1Q":"OOR","TUSDA2":-50,"TUSDA2Q":"OOR","TUSDA3":-50,"TUSDA3Q":"OOR","CmhVent":0,"CmhVentQ":null,"BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2025-05-12T10:01:32Z","PTIResult":"Failed","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":null,"ReeferAlarms":null}}}'

# Row 7: Sequence 6
$ws.Range("A7").Value = '2025-11-24T09:33:51.779Z'
$ws.Range("B7:I7").Value = $na
$ws.Range("J7").Value = '{"Sequence":6,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"17dfc11f-308c-4efd-8d45-00a9ec77313b","EventDtm":"2025-11-24T09:33:34Z","AppDtm":"2025-11-24T09:33:51Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-24T09:33:34Z","DeviceID":"JSGA623040295","GPSLockState":"LOCKED","GPSSatelliteCount":15,"GPSLastLock":0,"GPSLatitude":17.679088,"GPSLongitude":78.720172,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"49","LAC":"19328","CID":"254963202","RSSI":null,"ExtPower":true,"ExtPowerVoltage":29.6,"BatteryVoltage":8,"DeviceTemp":44,"RTDLOn":false,"VerFW":"W0206 1.91","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"CGMU2991560","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-24T09:33:34Z","AssetType":"Reefer","AssetID":"CGMU2991560","OEM":"CARRIER","TAmb":29.86,"TAmbQ":null,"TUSDA4":-50,"TUSDA4Q":"OOR","PctCO2":4,"PctCO2Q":null,"PctCO2Set":1,"PctCO2SetQ":null,"PSuc":-14.73,"PSucQ":"OOR","TDis":-181.25,"TDisQ":"OOR","FreqComp":null,"TSuc":0,"TSucQ":null,"MCond":1,"PCond":155.96,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":"Cool","SnCtrl":"04540798","AmpPhA":11.29,"AmpPhB":10.92,"AmpPhC":11.29,"PDis":-14.6,"PDisQ":"OOR","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":0,"PctExpVlvQ":null,"TEvap":5.46,"TEvapQ":null,"MCtrl3":null,"PctHtr":null,"PctHtrQ":null,"MEvapFanHS":1,"PctGasVlv":null,"PctGasVlvQ":null,"PctHum":98.81,"PctHumQ":null,"PctHumSet":0,"PctHumSetQ":"OOR","FreqLine":50,"FreqLineQ":null,"VLine1":410.67,"VLine2":null,"VLine3":null,"MEvapFanLS":0,"PctO2":-0.11,"PctO2Q":"OOR","PctO2Set":3,"PctO2SetQ":null,"MCtrl2":null,"TRtn1":5.02,"TRtn1Q":null,"TRtn2":5.01,"TRtn2Q":null,"TSet":4,"TSetQ":null,"VerSWMajor":"5180","VerSWMinor":null,"PctSucVlv":4.45,"PctSucVlvQ":null,"TSup1":3.89,"TSup1Q":null,"TSup2":3.91,"TSup2Q":null,"AmpTotalDraw":null,"AmpTotalDrawQ":null,"TUSDA1":-50,"TUSDA1Q":"OOR","TUSDA2":-50,"TUSDA2Q":"OOR","TUSDA3":-50,"TUSDA3Q":"OOR","CmhVent":0,"CmhVentQ":null,"BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2025-09-11T07:51:42Z","PTIResult":"Passed","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":null,"ReeferAlarms":[{"OemAlarm":53,"RCAlias":"AL53","RCSeverity":"Informational","Active":true}]}}}'

# Row 8: Sequence 7
$ws.Range("A8").Value = '2025-11-24T09:35:52.733Z'
$ws.Range("B8:I8").Value = $na
$ws.Range("J8").Value = '{"Sequence":7,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"0a8ac7fd-78ad-4734-b2b5-ba149232e16d","EventDtm":"2025-11-24T09:35:38Z","AppDtm":"2025-11-24T09:35:52Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-24T09:35:38Z","DeviceID":"JSGA622180045","GPSLockState":"LOCKED","GPSSatelliteCount":16,"GPSLastLock":0,"GPSLatitude":26.310566,"GPSLongitude":91.717636,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"405","MNC":"56","LAC":"7134","CID":"250551307","RSSI":null,"ExtPower":true,"ExtPowerVoltage":31.2,"BatteryVoltage":8.1,"DeviceTemp":37,"RTDLOn":false,"VerFW":"W0206 1.91","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"TDRU7151905","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-24T09:35:38Z","AssetType":"Reefer","AssetID":"TDRU7151905","OEM":"DAIKIN","TAmb":28.81,"TAmbQ":null,"TUSDA4":-46.4,"TUSDA4Q":"OOR","PctCO2":25.5,"PctCO2Q":"OOR","PctCO2Set":25.5,"PctCO2SetQ":"OOR","PSuc":-10,"PSucQ":"asProvided","TDis":73.3,"TDisQ":null,"FreqComp":null,"TSuc":-25.1,"TSucQ":null,"MCond":"On","PCond":null,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":null,"SnCtrl":null,"AmpPhA":null,"AmpPhB":null,"AmpPhC":null,"PDis":960,"PDisQ":"asProvided","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":72,"PctExpVlvQ":null,"TEvap":-0.06,"TEvapQ":null,"MCtrl3":"Modulation","PctHtr":null,"PctHtrQ":null,"MEvapFanHS":"On","PctGasVlv":0,"PctGasVlvQ":null,"PctHum":100.39,"PctHumQ":"OOR","PctHumSet":75,"PctHumSetQ":"configured","FreqLine":50,"FreqLineQ":null,"VLine1":393.8,"VLine2":null,"VLine3":null,"MEvapFanLS":"Off","PctO2":25.5,"PctO2Q":"OOR","PctO2Set":25.5,"PctO2SetQ":"OOR","MCtrl2":"Modulation","TRtn1":-0.06,"TRtn1Q":null,"TRtn2":5.4,"TRtn2Q":null,"TSet":4,"TSetQ":null,"VerSWMajor":"24C2","VerSWMinor":null,"PctSucVlv":92.07,"PctSucVlvQ":null,"TSup1":4.12,"TSup1Q":null,"TSup2":3.9,"TSup2Q":null,"AmpTotalDraw":11,"AmpTotalDrawQ":"asProvided","TUSDA1":-46.3,"TUSDA1Q":"OOR","TUSDA2":-46.3,"TUSDA2Q":"OOR","TUSDA3":-46.4,"TUSDA3Q":"OOR","CmhVent":1020,"CmhVentQ":"OOR","BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":null,"PTIResult":null,"TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":"UnitOff","ReeferAlarms":[{"OemAlarm":403,"RCAlias":"E403","RCSeverity":"Minor","Active":true},{"OemAlarm":409,"RCAlias":"E409","RCSeverity":"Informational","Active":true}]}}}'

# Row 9: Sequence 8
$ws.Range("A9").Value = '2025-11-24T09:43:13.486Z'
$ws.Range("B9:I9").Value = $na
$ws.Range("J9").Value = '{"Sequence":8,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"49130eb5-9e5c-4536-928f-5c3b18928cfa","EventDtm":"2025-11-24T09:42:59Z","AppDtm":"2025-11-24T09:43:13Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-24T09:42:59Z","DeviceID":"JSGA622340343","GPSLockState":"LOCKED","GPSSatelliteCount":17,"GPSLastLock":0,"GPSLatitude":12.80506,"GPSLongitude":77.662245,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"45","LAC":"9003","CID":"46904076","RSSI":null,"ExtPower":true,"ExtPowerVoltage":29.5,"BatteryVoltage":8,"DeviceTemp":33,"RTDLOn":false,"VerFW":"W0206 1.75","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"CCLU1035976","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-24T09:42:59Z","AssetType":"Reefer","AssetID":"CCLU1035976","OEM":"CARRIER","TAmb":32.86,"TAmbQ":null,"TUSDA4":-50,"TUSDA4Q":"OOR","PctCO2":5,"PctCO2Q":null,"PctCO2Set":1,"PctCO2SetQ":null,"PSuc":-14.73,"PSucQ":"OOR","TDis":-196.88,"TDisQ":"OOR","FreqComp":null,"TSuc":0,"TSucQ":null,"MCond":1,"PCond":103.09,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":"Cool","SnCtrl":"04904656","AmpPhA":11.25,"AmpPhB":10.79,"AmpPhC":11.25,"PDis":-14.6,"PDisQ":"OOR","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":0,"PctExpVlvQ":null,"TEvap":6.19,"TEvapQ":null,"MCtrl3":null,"PctHtr":null,"PctHtrQ":null,"MEvapFanHS":1,"PctGasVlv":null,"PctGasVlvQ":null,"PctHum":0,"PctHumQ":null,"PctHumSet":0,"PctHumSetQ":"OOR","FreqLine":50,"FreqLineQ":null,"VLine1":414.25,"VLine2":null,"VLine3":null,"MEvapFanLS":0,"PctO2":-0.06,"PctO2Q":"OOR","PctO2Set":3,"PctO2SetQ":null,"MCtrl2":null,"TRtn1":5.76,"TRtn1Q":null,"TRtn2":5.78,"TRtn2Q":null,"TSet":5,"TSetQ":null,"VerSWMajor":"5180","VerSWMinor":null,"PctSucVlv":4.84,"PctSucVlvQ":null,"TSup1":4.94,"TSup1Q":null,"TSup2":4.86,"TSup2Q":null,"AmpTotalDraw":null,"AmpTotalDrawQ":null,"TUSDA1":-50,"TUSDA1Q":"OOR","TUSDA2":-50,"TUSDA2Q":"OOR","TUSDA3":-50,"TUSDA3Q":"OOR","CmhVent":0,"CmhVentQ":null,"BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2025-05-12T10:01:32Z","PTIResult":"Failed","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":null,"ReeferAlarms":null}}}'

# Row 10: Sequence 9
$ws.Range("A10").Value = '2025-11-24T09:45:45.274Z'
$ws.Range("B10:I10").Value = $na
$ws.Range("J10").Value = '{"Sequence":9,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"4e06f62b-61d0-4b73-a9ef-a1c6df014d68","EventDtm":"2025-11-24T09:45:28Z","AppDtm":"2025-11-24T09:45:45Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-24T09:45:28Z","DeviceID":"JSGA623040278","GPSLockState":"LOCKED","GPSSatelliteCount":19,"GPSLastLock":0,"GPSLatitude":17.53317,"GPSLongitude":78.433499,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"49","LAC":"19315","CID":"231571714","RSSI":null,"ExtPower":true,"ExtPowerVoltage":32,"BatteryVoltage":8,"DeviceTemp":44,"RTDLOn":false,"VerFW":"W0206 1.91","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"CXRU1026520","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-24T09:45:28Z","AssetType":"Reefer","AssetID":"CXRU1026520","OEM":"DAIKIN","TAmb":31.81,"TAmbQ":null,"TUSDA4":-53.4,"TUSDA4Q":"OOR","PctCO2":25.5,"PctCO2Q":"OOR","PctCO2Set":25.5,"PctCO2SetQ":"OOR","PSuc":-10,"PSucQ":"asProvided","TDis":78,"TDisQ":null,"FreqComp":null,"TSuc":-27.8,"TSucQ":null,"MCond":"On","PCond":null,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":null,"SnCtrl":null,"AmpPhA":null,"AmpPhB":null,"AmpPhC":null,"PDis":940,"PDisQ":"asProvided","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":14,"PctExpVlvQ":null,"TEvap":3,"TEvapQ":null,"MCtrl3":"Modulation","PctHtr":null,"PctHtrQ":null,"MEvapFanHS":"On","PctGasVlv":0,"PctGasVlvQ":null,"PctHum":100,"PctHumQ":null,"PctHumSet":95,"PctHumSetQ":"off","FreqLine":null,"FreqLineQ":null,"VLine1":410.3,"VLine2":null,"VLine3":null,"MEvapFanLS":"Off","PctO2":25.5,"PctO2Q":"OOR","PctO2Set":25.5,"PctO2SetQ":"OOR","MCtrl2":"Modulation","TRtn1":6.31,"TRtn1Q":null,"TRtn2":6.4,"TRtn2Q":null,"TSet":5,"TSetQ":null,"VerSWMajor":"265A","VerSWMinor":null,"PctSucVlv":11.28,"PctSucVlvQ":null,"TSup1":5.12,"TSup1Q":null,"TSup2":5.1,"TSup2Q":null,"AmpTotalDraw":13,"AmpTotalDrawQ":"asProvided","TUSDA1":-53.4,"TUSDA1Q":"OOR","TUSDA2":-53.4,"TUSDA2Q":"OOR","TUSDA3":-53.4,"TUSDA3Q":"OOR","CmhVent":1020,"CmhVentQ":"error","BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2025-10-31T08:13:55Z","PTIResult":"Passed","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":"UnitOff","ReeferAlarms":null}}}'

# Row 11: Sequence 10
$ws.Range("A11").Value = '2025-11-24T09:48:38.612Z'
$ws.Range("B11:I11").Value = $na
$ws.Range("J11").Value = '{"Sequence":10,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"54ea56e8-0ba0-479f-b956-c087919f815f","EventDtm":"2025-11-24T09:48:24Z","AppDtm":"2025-11-24T09:48:38Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-24T09:48:24Z","DeviceID":"JSGA623040298","GPSLockState":"LOCKED","GPSSatelliteCount":20,"GPSLastLock":0,"GPSLatitude":17.657914,"GPSLongitude":83.101831,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"49","LAC":"31121","CID":"233290773","RSSI":"-75","ExtPower":true,"ExtPowerVoltage":28.6,"BatteryVoltage":8.1,"DeviceTemp":33,"RTDLOn":false,"VerFW":"W0206 1.91","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"KKFU6994964","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-24T09:48:24Z","AssetType":"Reefer","AssetID":"KKFU6994964","OEM":"CARRIER","TAmb":29.49,"TAmbQ":null,"TUSDA4":-50,"TUSDA4Q":"OOR","PctCO2":0,"PctCO2Q":null,"PctCO2Set":5,"PctCO2SetQ":null,"PSuc":16.22,"PSucQ":null,"TDis":-190.62,"TDisQ":"OOR","FreqComp":null,"TSuc":0,"TSucQ":null,"MCond":1,"PCond":226.04,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":"Cool","SnCtrl":"04783085","AmpPhA":17.41,"AmpPhB":17.77,"AmpPhC":17.41,"PDis":240.44,"PDisQ":null,"PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":0,"PctExpVlvQ":null,"TEvap":21.14,"TEvapQ":null,"MCtrl3":null,"PctHtr":null,"PctHtrQ":null,"MEvapFanHS":1,"PctGasVlv":null,"PctGasVlvQ":null,"PctHum":0.04,"PctHumQ":null,"PctHumSet":0,"PctHumSetQ":"OOR","FreqLine":50,"FreqLineQ":null,"VLine1":399.02,"VLine2":null,"VLine3":null,"MEvapFanLS":0,"PctO2":0,"PctO2Q":"asProvided","PctO2Set":10,"PctO2SetQ":null,"MCtrl2":null,"TRtn1":20.69,"TRtn1Q":null,"TRtn2":20.74,"TRtn2Q":null,"TSet":20,"TSetQ":null,"VerSWMajor":"5156","VerSWMinor":null,"PctSucVlv":22.02,"PctSucVlvQ":null,"TSup1":17.96,"TSup1Q":null,"TSup2":18.02,"TSup2Q":null,"AmpTotalDraw":null,"AmpTotalDrawQ":null,"TUSDA1":-50,"TUSDA1Q":"OOR","TUSDA2":-50,"TUSDA2Q":"OOR","TUSDA3":-50,"TUSDA3Q":"OOR","CmhVent":0,"CmhVentQ":null,"BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2024-11-29T10:30:33Z","PTIResult":"Passed","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":null,"ReeferAlarms":null}}}'

# Row 12: Sequence 11
$ws.Range("A12").Value = '2025-11-24T09:50:25.399Z'
$ws.Range("B12:I12").Value = $na
$ws.Range("J12").Value = '{"Sequence":11,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"5f7c45c9-f84d-4b13-a9ed-14334ddd0e7e","EventDtm":"2025-11-24T09:50:06Z","AppDtm":"2025-11-24T09:50:24Z","Events":["ACPowerOn"]},"DeviceData":{"DeviceDataDtm":"2025-11-24T09:50:06Z","DeviceID":"JSGA623040193","GPSLockState":"LOCKED","GPSSatelliteCount":13,"GPSLastLock":0,"GPSLatitude":28.678773,"GPSLongitude":77.59983,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"97","LAC":"1827","CID":"230077975","RSSI":"-71","ExtPower":true,"ExtPowerVoltage":29.1,"BatteryVoltage":8,"DeviceTemp":27,"RTDLOn":false,"VerFW":"W0206 1.75","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"SJKU4000017","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-24T09:50:06Z","AssetType":"Reefer","AssetID":"SJKU4000017","OEM":"CARRIER","TAmb":24.19,"TAmbQ":null,"TUSDA4":-50,"TUSDA4Q":"OOR","PctCO2":0,"PctCO2Q":null,"PctCO2Set":5,"PctCO2SetQ":null,"PSuc":70.74,"PSucQ":null,"TDis":33.09,"TDisQ":null,"FreqComp":null,"TSuc":20.87,"TSucQ":null,"MCond":0,"PCond":306.05,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":"Cool","SnCtrl":"04475864","AmpPhA":1.58,"AmpPhB":1.56,"AmpPhC":1.58,"PDis":84.31,"PDisQ":null,"PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":0,"PctExpVlvQ":null,"TEvap":21.14,"TEvapQ":null,"MCtrl3":null,"PctHtr":null,"PctHtrQ":null,"MEvapFanHS":1,"PctGasVlv":null,"PctGasVlvQ":null,"PctHum":61.62,"PctHumQ":null,"PctHumSet":0,"PctHumSetQ":"OOR","FreqLine":50,"FreqLineQ":null,"VLine1":409.54,"VLine2":null,"VLine3":null,"MEvapFanLS":0,"PctO2":0,"PctO2Q":"asProvided","PctO2Set":10,"PctO2SetQ":null,"MCtrl2":null,"TRtn1":21.69,"TRtn1Q":null,"TRtn2":21.67,"TRtn2Q":null,"TSet":17.2,"TSetQ":null,"VerSWMajor":"5370","VerSWMinor":null,"PctSucVlv":10,"PctSucVlvQ":null,"TSup1":20.51,"TSup1Q":null,"TSup2":20.41,"TSup2Q":null,"AmpTotalDraw":null,"AmpTotalDrawQ":null,"TUSDA1":-50,"TUSDA1Q":"OOR","TUSDA2":-50,"TUSDA2Q":"OOR","TUSDA3":-50,"TUSDA3Q":"OOR","CmhVent":0,"CmhVentQ":null,"BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2025-03-23T08:47:13Z","PTIResult":"Passed","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":null,"ReeferAlarms":null}}}'

# Row 13: Sequence 12
$ws.Range("A13").Value = '2025-11-24T09:51:17.401Z'
$ws.Range("B13:I13").Value = $na
$ws.Range("J13").Value = '{"Sequence":12,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"5f4b9f47-d565-42ee-9865-0558679c5965","EventDtm":"2025-11-24T09:48:34Z","AppDtm":"2025-11-24T09:51:17Z","Events":["DeviceIsStationary"]},"DeviceData":{"DeviceDataDtm":"2025-11-24T09:48:34Z","DeviceID":"JSGA622180064","GPSLockState":"LOCKED","GPSSatelliteCount":17,"GPSLastLock":0,"GPSLatitude":17.829046,"GPSLongitude":76.793505,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"405","MNC":"864","LAC":"91","CID":"6906161","RSSI":null,"ExtPower":false,"ExtPowerVoltage":0,"BatteryVoltage":7.4,"DeviceTemp":31,"RTDLOn":false,"VerFW":"W0206 1.91","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"CXRU1041571","LastAssetRunState":"Offline"},"ReeferData":null}}'
